# Auto-generated script applying scheduled market-price / profit recalculation updates
# to the per-job Leve profit tables (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 2091.5757
$ws.Range("I15").Value = 2091.5757
$ws.Range("K15").Value = 6274.7271
$ws.Range("M15").Value = -6105.7271
# Row 116
$ws.Range("H116").Value = 3103
$ws.Range("I116").Value = 2350.6
$ws.Range("J116").Value = 4839.3076
$ws.Range("K116").Value = 2350.6
$ws.Range("L116").Value = 4839.3076
$ws.Range("M116").Value = 1091.4
$ws.Range("N116").Value = -11723.3076
# Row 132
$ws.Range("H132").Value = 6901506.5
$ws.Range("I132").Value = 8004871.5
$ws.Range("J132").Value = 5476.5
$ws.Range("K132").Value = 24014614.5
$ws.Range("L132").Value = 16429.5
$ws.Range("M132").Value = -24012084.5
$ws.Range("N132").Value = -21489.5

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 19232922
$ws.Range("I2").Value = 31251374
$ws.Range("J2").Value = 3399.8
$ws.Range("K2").Value = 31251374
$ws.Range("L2").Value = 3399.8
$ws.Range("M2").Value = -31251261
$ws.Range("N2").Value = -3625.8
# Row 13
$ws.Range("H13").Value = 70004
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()
# Row 74
$ws.Range("H74").Value = 1103.7894
$ws.Range("I74").Value = 904.1429000000001
$ws.Range("J74").Value = 1662.8
$ws.Range("K74").Value = 904.1429000000001
$ws.Range("L74").Value = 1662.8
$ws.Range("M74").Value = -30.14290000000005
$ws.Range("N74").Value = -3410.8
# Row 77
$ws.Range("H77").Value = 1103.7894
$ws.Range("I77").Value = 904.1429000000001
$ws.Range("J77").Value = 1662.8
$ws.Range("K77").Value = 4520.7145
$ws.Range("L77").Value = 8314
$ws.Range("M77").Value = -152.7145
$ws.Range("N77").Value = -17050
# Row 116
$ws.Range("H116").Value = 19232922
$ws.Range("I116").Value = 31251374
$ws.Range("J116").Value = 3399.8
$ws.Range("K116").Value = 31251374
$ws.Range("L116").Value = 3399.8
$ws.Range("M116").Value = -31249080
$ws.Range("N116").Value = -7987.8
# Row 122
$ws.Range("H122").Value = 3145.9443
$ws.Range("I122").Value = 2070.3333
$ws.Range("K122").Value = 6210.999899999999
$ws.Range("M122").Value = -3760.999899999999

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 19232922
$ws.Range("I3").Value = 31251374
$ws.Range("J3").Value = 3399.8
$ws.Range("K3").Value = 31251374
$ws.Range("L3").Value = 3399.8
$ws.Range("M3").Value = -31251260
$ws.Range("N3").Value = -3627.8
# Row 86
$ws.Range("H86").Value = 1812.9445
$ws.Range("I86").Value = 1157.5
$ws.Range("K86").Value = 1157.5
$ws.Range("M86").Value = -34.5
# Row 89
$ws.Range("H89").Value = 1812.9445
$ws.Range("I89").Value = 1157.5
$ws.Range("K89").Value = 5787.5
$ws.Range("M89").Value = -171.5

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3115.8293
$ws.Range("I31").Value = 2217.8064
$ws.Range("J31").Value = 5899.7
$ws.Range("K31").Value = 2217.8064
$ws.Range("L31").Value = 5899.7
$ws.Range("M31").Value = -1922.8064
$ws.Range("N31").Value = -6489.7
# Row 34
$ws.Range("H34").Value = 3115.8293
$ws.Range("I34").Value = 2217.8064
$ws.Range("J34").Value = 5899.7
$ws.Range("K34").Value = 2217.8064
$ws.Range("L34").Value = 5899.7
$ws.Range("M34").Value = -2015.8064
$ws.Range("N34").Value = -6303.7
# Row 58
$ws.Range("H58").Value = 21741864
$ws.Range("I58").Value = 1332.2667
$ws.Range("J58").Value = 62505364
$ws.Range("K58").Value = 1332.2667
$ws.Range("L58").Value = 62505364
$ws.Range("M58").Value = -1129.2667
$ws.Range("N58").Value = -62505770
# Row 94
$ws.Range("H94").Value = 33335618
$ws.Range("I94").Value = 3749
$ws.Range("J94").Value = 38463596
$ws.Range("K94").Value = 3749
$ws.Range("L94").Value = 38463596
$ws.Range("M94").Value = -3298
$ws.Range("N94").Value = -38464498
# Row 136
$ws.Range("H136").Value = 21741864
$ws.Range("I136").Value = 1332.2667
$ws.Range("J136").Value = 62505364
$ws.Range("K136").Value = 3996.800099999999
$ws.Range("L136").Value = 187516092
$ws.Range("M136").Value = -1446.800099999999
$ws.Range("N136").Value = -187521192

$ws = $wb.Worksheets.Item("CUL")
# Row 87
$ws.Range("H87").Value = 9979.333000000001
$ws.Range("I87").Value = 7469
$ws.Range("K87").Value = 22407
$ws.Range("M87").Value = -21159
# Row 90
$ws.Range("H90").Value = 9979.333000000001
$ws.Range("I90").Value = 7469
$ws.Range("K90").Value = 67221
$ws.Range("M90").Value = -60981
# Row 129
$ws.Range("H129").Value = 1974.5714
$ws.Range("I129").Value = 1744
$ws.Range("J129").Value = 2046.625
$ws.Range("K129").Value = 5232
$ws.Range("L129").Value = 6139.875
$ws.Range("M129").Value = -232
$ws.Range("N129").Value = -16139.875
# Row 131
$ws.Range("H131").Value = 1037.8214
$ws.Range("I131").Value = 544
$ws.Range("J131").Value = 1312.1666
$ws.Range("K131").Value = 1632
$ws.Range("L131").Value = 3936.4998
$ws.Range("M131").Value = 3408
$ws.Range("N131").Value = -14016.4998

$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 4284.36
$ws.Range("I122").Value = 4364.9287
$ws.Range("J122").Value = 4181.8184
$ws.Range("K122").Value = 13094.7861
$ws.Range("L122").Value = 12545.4552
$ws.Range("M122").Value = -10644.7861
$ws.Range("N122").Value = -17445.4552
# Row 126
$ws.Range("H126").Value = 2320.487
$ws.Range("I126").Value = 1361.5
$ws.Range("K126").Value = 4084.5
$ws.Range("M126").Value = -1614.5

$ws = $wb.Worksheets.Item("LTW")
# Row 3
$ws.Range("H3").Value = 70005
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 70005
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 70005
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = -70229
# Row 11
$ws.Range("H11").Value = 57506.75
$ws.Range("I11").Value = 20006
$ws.Range("J11").Value = 70007
$ws.Range("K11").Value = 20006
$ws.Range("L11").Value = 70007
$ws.Range("M11").Value = -19866
$ws.Range("N11").Value = -70287
# Row 15
$ws.Range("H15").Value = 70005
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 70005
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 70005
$ws.Range("M15").ClearContents()
$ws.Range("N15").Value = -70345
# Row 17
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()
# Row 21
$ws.Range("H21").Value = 47006.668
$ws.Range("I21").Value = 1006
$ws.Range("K21").Value = 1006
$ws.Range("M21").Value = -832
# Row 69
$ws.Range("H69").Value = 40000
$ws.Range("J69").Value = 40000
$ws.Range("L69").Value = 40000
$ws.Range("N69").Value = -41622
# Row 72
$ws.Range("H72").Value = 40000
$ws.Range("J72").Value = 40000
$ws.Range("L72").Value = 120000
$ws.Range("N72").Value = -128112
# Row 132
$ws.Range("H132").Value = 2869.9119
$ws.Range("I132").Value = 1632.1666
$ws.Range("J132").Value = 4262.375
$ws.Range("K132").Value = 4896.4998
$ws.Range("L132").Value = 12787.125
$ws.Range("M132").Value = -2366.4998
$ws.Range("N132").Value = -17847.125

$ws = $wb.Worksheets.Item("WVR")
# Row 7
$ws.Range("H7").Value = 42122.8
$ws.Range("I7").Value = 302
$ws.Range("J7").Value = 70003.336
$ws.Range("K7").Value = 302
$ws.Range("L7").Value = 70003.336
$ws.Range("M7").Value = -189
$ws.Range("N7").Value = -70229.336
# Row 11
$ws.Range("H11").Value = 57003
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 57003
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 57003
$ws.Range("M11").ClearContents()
$ws.Range("N11").Value = -57287
# Row 12
$ws.Range("H12").Value = 56965.4
$ws.Range("I12").Value = 4806
$ws.Range("J12").Value = 70005.25
$ws.Range("K12").Value = 4806
$ws.Range("L12").Value = 70005.25
$ws.Range("M12").Value = -4664
$ws.Range("N12").Value = -70289.25
# Row 13
$ws.Range("H13").Value = 69633.336
$ws.Range("I13").Value = 38888
$ws.Range("K13").Value = 38888
$ws.Range("M13").Value = -38748
# Row 132
$ws.Range("H132").Value = 15441.303
$ws.Range("I132").Value = 3955.7666
$ws.Range("J132").Value = 41946.383
$ws.Range("K132").Value = 11867.2998
$ws.Range("L132").Value = 125839.149
$ws.Range("M132").Value = -9337.299800000001
$ws.Range("N132").Value = -130899.149
